$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Music Subscription -> rent
$ws.Range("A2").Value = "rent"
$ws.Range("B2").Value = 20000
$ws.Range("C2").Value = 45778.22928240741

# Update row 3: Groceries -> rent
$ws.Range("A3").Value = "rent"
$ws.Range("B3").Value = 20000
$ws.Range("C3").Value = 45778.22928240741

# Update row 4: Rent -> lunch
$ws.Range("A4").Value = "lunch"
$ws.Range("B4").Value = 3000
$ws.Range("C4").Value = 45778.22928240741

# Remove row 5 entirely (was spoon/600)
$ws.Range("A5:C5").Delete()
